$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new values
$ws.Range("A2").Value = "t13"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "Computer Science and Engineering"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Delete row 3 entirely (the "Cool Kids Club" row)
$ws.Rows("3").Delete()
